# ------------------------------------------------------------------
# Commit: "added code for table 4 (outcome by etiology) and added
#          initial treatment to table 3"
#
# 1. Demographics: reorder the "pacific islander" line within the
#    Ethnicity cells for columns D and F.
# 2. Demographics: apply a centered-header look to the label/header
#    cells (this is what introduces the extra centered cell style in
#    the workbook's style table).
# 3. Outcome: insert a new "Initial Treatment" row (pushing "Final
#    Treatment" and "Outcome" down one row).
# 4. Add a brand-new "Outcome by Etiology" worksheet (after "Outcome")
#    with Initial Treatment / Final Treatment / Outcome broken out by
#    CSA etiology instead of by CSA severity.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ==================================================================
# 1) Demographics: fix the ordering of "pacific islander" in the
#    Ethnicity row for the two CSA-severity columns.
# ==================================================================
$demo = $wb.Worksheets.Item("Demographics")

$demo.Range("D5").Value = "not hispanic/latino = 40 (48.8%)`nwhite = 35 (42.7%)`nhispanic = 4 (4.9%)`nasian = 1 (1.2%)`nnative american = 1 (1.2%)`npacific islander = 1 (1.2%)"
$demo.Range("F5").Value = "not hispanic/latino = 104 (58.4%)`nwhite = 58 (32.6%)`nhispanic = 11 (6.2%)`npacific islander = 2 (1.1%)`nasian = 2 (1.1%)`nnative american = 1 (0.6%)"

# ==================================================================
# 2) Demographics: center the header/label cells (adds the new
#    centered style to the workbook; Etiology/Outcome/new sheet keep
#    their original left/top-aligned header look).
# ==================================================================
$demo.Range("A1:F1,A2:A9").HorizontalAlignment = -4108

# ==================================================================
# 3) Outcome: insert "Initial Treatment" as the new row 3.
# ==================================================================
$outcome = $wb.Worksheets.Item("Outcome")

$outcome.Rows.Item(3).Insert()

# Give the new A3 label cell the same look as the other label cells
# in column A (border + bold + wrap + vertical-center).
$outcome.Range("A2").Copy($outcome.Range("A3"))

$row3 = @{}
$row3["A3"] = "Initial Treatment"
$row3["B3"] = "cpap = 453 (73.4%)`nasv = 109 (17.7%)`nunknown = 26 (4.2%)`nnone = 15 (2.4%)`nO2 = 13 (2.1%)`nother = 1 (0.2%)`nbipap = 0 (0.0%)`nivaps = 0 (0.0%)"
$row3["C3"] = "cpap = 9 (40.9%)`nasv = 8 (36.4%)`nunknown = 2 (9.1%)`nO2 = 2 (9.1%)`nnone = 1 (4.5%)`nother = 0 (0.0%)`nbipap = 0 (0.0%)`nivaps = 0 (0.0%)"
$row3["D3"] = "cpap = 45 (54.9%)`nasv = 25 (30.5%)`nO2 = 7 (8.5%)`nnone = 3 (3.7%)`nunknown = 2 (2.4%)`nother = 0 (0.0%)`nbipap = 0 (0.0%)`nivaps = 0 (0.0%)"
$row3["E3"] = "cpap = 250 (74.6%)`nasv = 62 (18.5%)`nunknown = 11 (3.3%)`nnone = 9 (2.7%)`nO2 = 3 (0.9%)`nother = 0 (0.0%)`nbipap = 0 (0.0%)`nivaps = 0 (0.0%)"
$row3["F3"] = "cpap = 149 (83.7%)`nasv = 14 (7.9%)`nunknown = 11 (6.2%)`nnone = 2 (1.1%)`nother = 1 (0.6%)`nO2 = 1 (0.6%)`nbipap = 0 (0.0%)`nivaps = 0 (0.0%)"

foreach ($ref in $row3.Keys) {
    $outcome.Range($ref).Value = $row3[$ref]
}

# ==================================================================
# 4) Add the new "Outcome by Etiology" worksheet at the end of the
#    workbook and populate it.
# ==================================================================
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$etio2 = $wb.Worksheets.Add($null, $lastSheet)
$etio2.Name = "Outcome by Etiology"

# Seed formatting by copying the (now 5-row) Outcome sheet's label
# column + header row formats across, then overwrite the values.
$outcome.Range("A1:F2").Copy($etio2.Range("A1:F2"))
$outcome.Range("F1").Copy($etio2.Range("G1:H1"))
$outcome.Range("A3").Copy($etio2.Range("A3"))
$outcome.Range("A4").Copy($etio2.Range("A4"))
$outcome.Range("A5").Copy($etio2.Range("A5"))
$etio2.Range("B2:F2").ClearContents()

$headers = @{}
$headers["B1"] = "All, n=617"
$headers["C1"] = "Neurologic Contributor, n=98"
$headers["D1"] = "Cardiac Contributor, n=136"
$headers["E1"] = "Medication Contributor, n=67"
$headers["F1"] = "Treatment Emergent, n=239"
$headers["G1"] = "OSA-associated Centrals, n=105"
$headers["H1"] = "Primary CSA, n=13"
foreach ($ref in $headers.Keys) {
    $etio2.Range($ref).Value = $headers[$ref]
}

$etio2.Range("A3").Value = "Initial Treatment"
$etio2.Range("A4").Value = "Final Treatment"
$etio2.Range("A5").Value = "Outcome"

$vals = @{}
$vals["B3"] = "cpap = 453 (73.4%)`nasv = 109 (17.7%)`nunknown = 26 (4.2%)`nnone = 15 (2.4%)`nO2 = 13 (2.1%)`nother = 1 (0.2%)`nbipap = 0 (0.0%)`nivaps = 0 (0.0%)"
$vals["B4"] = "cpap = 312 (50.6%)`nasv = 205 (33.2%)`nbipap = 51 (8.3%)`nnone = 23 (3.7%)`nO2 = 18 (2.9%)`nother = 7 (1.1%)`nivaps = 1 (0.2%)"
$vals["B5"] = "resolved w/ cpap = 245 (39.7%)`nfailed cpap = 141 (22.9%)`nn/a = 127 (20.6%)`nnon-compliant = 88 (14.3%)`nnever started on cpap = 14 (2.3%)`nresolved w/bipap = 2 (0.3%)"
$vals["C3"] = "cpap = 67 (68.4%)`nasv = 24 (24.5%)`nunknown = 6 (6.1%)`nnone = 1 (1.0%)`nother = 0 (0.0%)`nO2 = 0 (0.0%)`nbipap = 0 (0.0%)`nivaps = 0 (0.0%)"
$vals["C4"] = "asv = 44 (44.9%)`ncpap = 43 (43.9%)`nbipap = 6 (6.1%)`nnone = 2 (2.0%)`nother = 1 (1.0%)`nO2 = 1 (1.0%)`nivaps = 1 (1.0%)"
$vals["C5"] = "resolved w/ cpap = 28 (28.6%)`nfailed cpap = 26 (26.5%)`nn/a = 25 (25.5%)`nnon-compliant = 18 (18.4%)`nnever started on cpap = 1 (1.0%)`nresolved w/bipap = 0 (0.0%)"
$vals["D3"] = "cpap = 99 (72.8%)`nasv = 20 (14.7%)`nO2 = 7 (5.1%)`nunknown = 6 (4.4%)`nnone = 4 (2.9%)`nother = 0 (0.0%)`nbipap = 0 (0.0%)`nivaps = 0 (0.0%)"
$vals["D4"] = "cpap = 63 (46.3%)`nasv = 44 (32.4%)`nbipap = 14 (10.3%)`nO2 = 10 (7.4%)`nnone = 5 (3.7%)`nother = 0 (0.0%)`nivaps = 0 (0.0%)"
$vals["D5"] = "resolved w/ cpap = 46 (33.8%)`nfailed cpap = 40 (29.4%)`nn/a = 29 (21.3%)`nnon-compliant = 16 (11.8%)`nnever started on cpap = 5 (3.7%)`nresolved w/bipap = 0 (0.0%)"
$vals["E3"] = "cpap = 36 (53.7%)`nasv = 22 (32.8%)`nunknown = 4 (6.0%)`nO2 = 3 (4.5%)`nnone = 2 (3.0%)`nother = 0 (0.0%)`nbipap = 0 (0.0%)`nivaps = 0 (0.0%)"
$vals["E4"] = "asv = 40 (59.7%)`ncpap = 15 (22.4%)`nO2 = 4 (6.0%)`nnone = 3 (4.5%)`nbipap = 3 (4.5%)`nother = 1 (1.5%)`nivaps = 1 (1.5%)"
$vals["E5"] = "n/a = 22 (32.8%)`nfailed cpap = 21 (31.3%)`nnon-compliant = 10 (14.9%)`nresolved w/ cpap = 9 (13.4%)`nnever started on cpap = 5 (7.5%)`nresolved w/bipap = 0 (0.0%)"
$vals["F3"] = "cpap = 188 (78.7%)`nasv = 31 (13.0%)`nunknown = 11 (4.6%)`nnone = 5 (2.1%)`nO2 = 3 (1.3%)`nother = 1 (0.4%)`nbipap = 0 (0.0%)`nivaps = 0 (0.0%)"
$vals["F4"] = "cpap = 136 (56.9%)`nasv = 62 (25.9%)`nbipap = 24 (10.0%)`nnone = 9 (3.8%)`nother = 4 (1.7%)`nO2 = 4 (1.7%)`nivaps = 0 (0.0%)"
$vals["F5"] = "resolved w/ cpap = 119 (49.8%)`nfailed cpap = 48 (20.1%)`nn/a = 39 (16.3%)`nnon-compliant = 31 (13.0%)`nresolved w/bipap = 2 (0.8%)`nnever started on cpap = 0 (0.0%)"
$vals["G3"] = "cpap = 81 (77.1%)`nasv = 20 (19.0%)`nnone = 3 (2.9%)`nunknown = 1 (1.0%)`nother = 0 (0.0%)`nO2 = 0 (0.0%)`nbipap = 0 (0.0%)`nivaps = 0 (0.0%)"
$vals["G4"] = "cpap = 65 (61.9%)`nasv = 31 (29.5%)`nnone = 4 (3.8%)`nbipap = 4 (3.8%)`nother = 1 (1.0%)`nO2 = 0 (0.0%)`nivaps = 0 (0.0%)"
$vals["G5"] = "resolved w/ cpap = 51 (48.6%)`nn/a = 21 (20.0%)`nnon-compliant = 16 (15.2%)`nfailed cpap = 15 (14.3%)`nnever started on cpap = 2 (1.9%)`nresolved w/bipap = 0 (0.0%)"
$vals["H3"] = "cpap = 8 (61.5%)`nasv = 3 (23.1%)`nunknown = 2 (15.4%)`nother = 0 (0.0%)`nnone = 0 (0.0%)`nO2 = 0 (0.0%)`nbipap = 0 (0.0%)`nivaps = 0 (0.0%)"
$vals["H4"] = "cpap = 6 (46.2%)`nasv = 6 (46.2%)`nbipap = 1 (7.7%)`nother = 0 (0.0%)`nnone = 0 (0.0%)`nO2 = 0 (0.0%)`nivaps = 0 (0.0%)"
$vals["H5"] = "resolved w/ cpap = 5 (38.5%)`nnon-compliant = 2 (15.4%)`nnever started on cpap = 2 (15.4%)`nn/a = 2 (15.4%)`nfailed cpap = 2 (15.4%)`nresolved w/bipap = 0 (0.0%)"

foreach ($ref in $vals.Keys) {
    $etio2.Range($ref).Value = $vals[$ref]
}
